$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.324.56'
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.581.56'
$ws.Range("E3").Value = '  -2.17%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.21'
$ws.Range("E5").Value = '  -2.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.42'
$ws.Range("E6").Value = '  +2.17%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -0.37%  '
$ws.Range("E9").Value = '  +0.47%  '
$ws.Range("E10").Value = '  +1.90%  '
$ws.Range("E11").Value = '  +0.48%  '
$ws.Range("E12").Value = '  -0.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.60'
$ws.Range("E13").Value = '  +0.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.044.76'
$ws.Range("E14").Value = '  -2.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '63.149.49'
$ws.Range("E15").Value = '  -0.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000156'
$ws.Range("E16").Value = '  +5.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.606.41'
$ws.Range("E17").Value = '  -1.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.23'
$ws.Range("E18").Value = '  +4.16%  '
$ws.Range("E19").Value = '  +4.16%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '345.69'
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("E21").Value = '  -0.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.31'
$ws.Range("E23").Value = '  +1.55%  '
$ws.Range("E24").Value = '  +1.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.27'
$ws.Range("E25").Value = '  +0.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.68'
$ws.Range("E26").Value = '  -0.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '564.56'
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.06'
$ws.Range("E28").Value = '  -0.11%  '
$ws.Range("E29").Value = '  +0.35%  '
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("E31").Value = '  -1.11%  '
$ws.Range("E32").Value = '  -0.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.22'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '166.74'
$ws.Range("E35").Value = '  -1.68%  '
$ws.Range("E36").Value = '  +1.94%  '
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.48'
$ws.Range("E38").Value = '  +1.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.92'
$ws.Range("E39").Value = '  -0.81%  '
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '166.19'
$ws.Range("E41").Value = '  +0.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.54'
$ws.Range("E42").Value = '  -1.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.97'
$ws.Range("E43").Value = '  +4.93%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.79'
$ws.Range("E44").Value = '  +4.30%  '
$ws.Range("E45").Value = '  +2.97%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.07'
$ws.Range("E46").Value = '  +3.21%  '
$ws.Range("E47").Value = '  +0.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0252'
$ws.Range("E48").Value = '  +3.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0961'
$ws.Range("E49").Value = '  +0.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.14'
$ws.Range("E50").Value = '  +2.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0231'
$ws.Range("E51").Value = '  +17.21%  '
